# Add season record columns (Wins, Losses, Ties) to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells: AD1 = "Wins", AE1 = "Losses", AF1 = "Ties"
# Apply the same look as the rest of the header row (bold, centered,
# top-aligned, thin border on all sides).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Every data row (2 through 48) gets the same season record: 89 wins,
# 73 losses, 0 ties. (Columns AD=30, AE=31, AF=32)
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 89
    $ws.Cells.Item($row, 31).Value = 73
    $ws.Cells.Item($row, 32).Value = 0
}
